# Updates cryptos list figures (price/volume columns) per the
# scraped-data refresh. Cells keep their original text cell type;
# for Price-column values that parse as plain numbers we pin the
# cell to Text format first so Excel does not silently convert the
# literal (e.g. "58.84") into a floating-point number on write, then
# restore the cell style so no stray formatting diff is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '38.005.46'
$ws.Range('E2').Value = '  +2.20%  '
$ws.Range('D3').Value = '2.052.66'
$ws.Range('E3').Value = '  +1.22%  '
$ws.Range('E4').Value = '  -0.32%  '
$ws.Range('D5').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '229.53'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.14%  '
$ws.Range('E6').Value = '  +1.18%  '
$ws.Range('D7').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.84'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +6.81%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.385'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.75%  '
$ws.Range('D10').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0808'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.76%  '
$ws.Range('E11').Value = '  +1.47%  '
$ws.Range('D12').Value = '2.356.11'
$ws.Range('E12').Value = '  +1.71%  '
$ws.Range('D13').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.66'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.96%  '
$ws.Range('D14').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.88'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.14%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.29'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.81%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.752'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.18%  '
$ws.Range('D17').Value = '2.050.25'
$ws.Range('E17').Value = '  +0.01%  '
$ws.Range('D18').Value = '37.939.62'
$ws.Range('E19').Value = '  -3.52%  '
$ws.Range('D20').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.69'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.22%  '
$ws.Range('D21').Value = '0.0₃0834'
$ws.Range('E21').Value = '  +2.01%  '
$ws.Range('D22').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '224.61'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.34%  '
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('E25').Value = '  +2.59%  '
$ws.Range('E26').Value = '  +0.90%  '
$ws.Range('D27').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '166.24'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.33%  '
$ws.Range('D28').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.133'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.18%  '
$ws.Range('D29').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.03'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.53%  '
$ws.Range('E30').Value = '  +0.56%  '
$ws.Range('E31').Value = '  +1.32%  '
$ws.Range('D32').Style = 'Normal'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.52'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.13%  '
$ws.Range('D33').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.59'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.23%  '
$ws.Range('E34').Value = '  +10.65%  '
$ws.Range('E35').Value = '  -0.51%  '
$ws.Range('D36').Style = 'Normal'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.34'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.63%  '
$ws.Range('E37').Value = '  +9.72%  '
$ws.Range('E38').Value = '  +5.55%  '
$ws.Range('E39').Value = '  +0.14%  '
$ws.Range('D40').Value = '1.487.57'
$ws.Range('E40').Value = '  +1.22%  '
$ws.Range('E41').Value = '  +0.75%  '
$ws.Range('E42').Value = '  +1.30%  '
$ws.Range('D43').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.86'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.58%  '
$ws.Range('D44').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.55'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.61%  '
$ws.Range('E45').Value = '  +1.56%  '
$ws.Range('E46').Value = '  -1.01%  '
$ws.Range('D47').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.15'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +14.07%  '
$ws.Range('E48').Value = '  +0.45%  '
$ws.Range('E49').Value = '  +1.04%  '
$ws.Range('E50').Value = '  -2.74%  '
$ws.Range('D51').Value = '2.244.33'
$ws.Range('E51').Value = '  +1.51%  '
